$d = $word.ActiveDocument

# 1. Split the merged "16h20 – 16h30" run after "16h20" is actually a
#    separate paragraph already containing "16h20 – 16h30". The diff
#    instead inserts a new " " run and a new "16h20" run right after the
#    "13h15 – " run (i.e. duplicates the start time onto the first line).
$d.Content.Find.Execute("13h15 – ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "13h15 –  16h20", 2)

# 2. Shift the run-split point in "Account Takeover" by one character:
#    "Account Takeov" + "er"  ->  "Account Takeove" + "r"
$d.Content.Find.Execute("Account Takeover", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Account Takeover", 2)
